$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix the timestamp on the last existing block (rows 226-239) ---
for ($r = 226; $r -le 239; $r++) {
    $ws.Cells.Item($r, 4).Value = 44232.00822890046
}

# --- 2) Append a new 14-row block (rows 240-253) following the same
#         repeating pattern (Nombre / URL / Disponibilidad / Fecha) ---

# Column A text (Nombre) per new row, 240..253
$names = @{
    240 = "Odoo"
    241 = "Blackbox"
    242 = "PowerBI"
    243 = "Dropbox"
    244 = "Odoo"
    245 = "GEE"
    246 = "UtilidadesOdoo"
    247 = "Filtros Dashboard"
    248 = "MapStore"
    249 = "GeoServer"
    250 = "Tomcat"
    251 = "Shiny"
    252 = "Github"
    253 = "EZ Exporter"
}

# Column B text (URL shown in the cell) per new row
$urlText = @{
    240 = "https://www.dataintelligence-group.com/"
    241 = "https://serviciodashboard.azurewebsites.net/"
    242 = "https://powerbi.microsoft.com/es-es/"
    243 = "https://www.dropbox.com/"
    244 = "https://dataintelligence.store/"
    245 = "https://app-data-i.users.earthengine.app/"
    246 = "https://odooutil.azurewebsites.net/"
    247 = "https://filtradordashboard.azurewebsites.net/"
    248 = "https://ide.dataintelligence-group.com/mapstore/#/"
    249 = "https://ide.dataintelligence-group.com/geoserver/web/?0"
    250 = "https://ide.dataintelligence-group.com/"
    251 = "https://rpubs.com/dataintelligence/"
    252 = "https://github.com/Sud-Austral/"
    253 = "https://ezexporter.highviewapps.com/exports/export-profile/"
}

# Hyperlink target address (the relationship Target, before any '#sub')
$hlTarget = @{
    240 = "https://www.dataintelligence-group.com/"
    241 = "https://serviciodashboard.azurewebsites.net/"
    242 = "https://powerbi.microsoft.com/es-es/"
    243 = "https://www.dropbox.com/"
    244 = "https://dataintelligence.store/"
    245 = "https://app-data-i.users.earthengine.app/"
    246 = "https://odooutil.azurewebsites.net/"
    247 = "https://filtradordashboard.azurewebsites.net/"
    248 = "https://ide.dataintelligence-group.com/mapstore/"
    249 = "https://ide.dataintelligence-group.com/geoserver/web/?0"
    250 = "https://ide.dataintelligence-group.com/"
    251 = "https://rpubs.com/dataintelligence/"
    252 = "https://github.com/Sud-Austral/"
    253 = "https://ezexporter.highviewapps.com/exports/export-profile/"
}

# Hyperlink sub-address ("location"), only the MapStore row (248) has one
$hlSubAddress = @{
    248 = "/"
}

$newDate = 44232.02931486353

for ($r = 240; $r -le 253; $r++) {
    $ws.Cells.Item($r, 1).Value = $names[$r]
    $ws.Cells.Item($r, 2).Value = $urlText[$r]
    $ws.Cells.Item($r, 2).Style = "Hyperlink"
    $ws.Cells.Item($r, 3).Value = "Disponible"
    $ws.Cells.Item($r, 4).Value = $newDate
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    if ($hlSubAddress.ContainsKey($r)) {
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 2), $hlTarget[$r], $hlSubAddress[$r])
    } else {
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 2), $hlTarget[$r])
    }
}

# Adding a hyperlink re-applies its own "Hyperlink" cell style, which lands
# on a freshly minted (but equivalent) style record instead of reusing the
# one already used by every pre-existing hyperlink cell (style index 2).
# Re-apply the named style across the whole hyperlink column so the new
# cells fold back onto that same, already-existing style.
$ws.Range("B2:B253").Style = "Hyperlink"
